# Insert two new data rows at the top of the data block (rows 35-36),
# pushing the existing rows 35..144 down to 37..146.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35:A36").EntireRow.Insert()

# --- New row 35 ---
$ws.Range("A35").Value = 10
$ws.Range("B35").Value = "Vega Modelo de Temuco"
$ws.Range("C35").Value = "La Araucanía"
$ws.Range("D35").Value = 45246
$ws.Range("E35").Value = 9
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100101
$ws.Range("H35").Value = "Berries"
$ws.Range("I35").Value = 100101001
$ws.Range("J35").Value = "Arándano (blue)"
$ws.Range("K35").Value = "Sin especificar"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 100
$ws.Range("N35").Value = 8000
$ws.Range("O35").Value = 8000
$ws.Range("P35").Value = 8000
$ws.Range("Q35").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R35").Value = "Provincia de Quillota"
$ws.Range("S35").Value = 5333
$ws.Range("T35").Value = 1.5

# --- New row 36 ---
$ws.Range("A36").Value = 10
$ws.Range("B36").Value = "Vega Modelo de Temuco"
$ws.Range("C36").Value = "La Araucanía"
$ws.Range("D36").Value = 45246
$ws.Range("E36").Value = 9
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100101
$ws.Range("H36").Value = "Berries"
$ws.Range("I36").Value = 100101001
$ws.Range("J36").Value = "Arándano (blue)"
$ws.Range("K36").Value = "Sin especificar"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 650
$ws.Range("N36").Value = 4800
$ws.Range("O36").Value = 5000
$ws.Range("P36").Value = 4923
$ws.Range("Q36").Value = "$/kilo"
$ws.Range("R36").Value = "Región del Maule"
$ws.Range("S36").Value = 4923
$ws.Range("T36").Value = 1
